# Actualizada grafica burn up
# Updates the "Hoja3" progress-tracking sheet: several tasks that were
# previously marked "No" (not yet at revision 3) are now marked complete
# (value 3), and the cumulative "Tareas completadas" count for the
# Revision-3 milestone (D8) is filled in with its COUNTIF formula, which
# in turn feeds the burn-up chart on Hoja1.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja3")

# Mark these tasks' current progress (column E) as having reached
# revision 3 (previously they held the "No" shared-string marker).
$ws.Range("E24").Value = 3
$ws.Range("E28").Value = 3
$ws.Range("E29").Value = 3
$ws.Range("E30").Value = 3
$ws.Range("E31").Value = 3
$ws.Range("E39").Value = 3

# Fill in the cumulative "completed" count for the Revision-3 row, mirroring
# the existing formula pattern used by the other milestone rows / by column E.
$ws.Range("D8").Formula = "=COUNTIF(E16:E54,""=3"")+D7"

# Update the last-saved cell selection on Hoja3 to reflect where the author
# left off, then restore the originally-active sheet (Hoja1) as the active
# tab so only Hoja3's own stored selection changes.
$ws.Range("D9").Select()
$wb.Worksheets.Item("Hoja1").Activate()

$wb.Application.Calculate()
